$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing data down
$ws.Rows.Item(1).Insert()

# Add header row
$ws.Range("A1").Value = "STATE"
$ws.Range("B1").Value = "COUNT"

# Restore the view's selected cell / scroll position to match the saved state
[void]$ws.Range("F14").Select()
